$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column N ("FinalTest1" / "int32[]") - filled first across the data rows
$ws.Range("N1").Value = "FinalTest1"
$ws.Range("N2").Value = "int32[]"
$ws.Range("N4").Value = "12|13|14"
$ws.Range("N5").Value = "31|33|34"
$ws.Range("N6").Value = "31|33|34"

# Column O ("FinalTest2" / "<string, string>")
$ws.Range("O1").Value = "FinalTest2"
$ws.Range("O2").Value = "<string, string>"

# Column P ("FinalTest3" / "dic<string, string>")
$ws.Range("P1").Value = "FinalTest3"
$ws.Range("P2").Value = "dic<string, string>"

# Column Q ("FinalTest4" / "int32")
$ws.Range("Q1").Value = "FinalTest4"
$ws.Range("Q2").Value = "int32"

# Row 3 comments for all four new columns
$ws.Range("N3").Value = "int32[]类型测试"
$ws.Range("O3").Value = "<string, string>类型测试"
$ws.Range("P3").Value = "dic<string, string>类型测试"
$ws.Range("Q3").Value = "int32类型测试"

# Remaining data cells - reuse existing shared string "s:aa|d:bb|f:cc" and plain numbers
$ws.Range("O4").Value = "s:aa|d:bb|f:cc"
$ws.Range("P4").Value = "s:aa|d:bb|f:cc"
$ws.Range("Q4").Value = 55

$ws.Range("O5").Value = "s:aa|d:bb|f:cc"
$ws.Range("P5").Value = "s:aa|d:bb|f:cc"
$ws.Range("Q5").Value = 66

$ws.Range("O6").Value = "s:aa|d:bb|f:cc"
$ws.Range("P6").Value = "s:aa|d:bb|f:cc"
$ws.Range("Q6").Value = 777

# Update selection to match the target workbook state
$null = $ws.Range("U5").Select()
